$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# New matchup data rows (Player_1, Points_1, Player_2, Points_2) for rows 912-926
$data = @(
    @(4, 12, 2, 8),
    @(5, 7, 4, 13),
    @(4, 12, 6, 8),
    @(5, 13, 9, 7),
    @(7, 6, 5, 14),
    @(5, 8, 4, 12),
    @(6, 19, 5, 1),
    @(2, 19, 1, 1),
    @(6, 4, 5, 16),
    @(4, 6, 8, 14),
    @(5, 7, 4, 13),
    @(4, 16, 5, 4),
    @(4, 12, 3, 8),
    @(7, 6, 6, 14),
    @(3, 18, 4, 2)
)

$startRow = 912
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Update the view: scroll position and selection, matching saved workbook state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 911
$ws.Range("F921:G923").Select()
